$wb = $excel.ActiveWorkbook

# The two sheets that carry this event table ("展览" and "全部类型") are
# identical in content and receive identical edits.
$targetSheets = @(1, 4)

foreach ($sheetIndex in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # --- 1. Updated "想去人数" (F column) counts for existing rows ---
    $ws.Cells.Item(2, 6).Value = 1927
    $ws.Cells.Item(4, 6).Value = 114
    $ws.Cells.Item(7, 6).Value = 1606
    $ws.Cells.Item(9, 6).Value = 635
    $ws.Cells.Item(13, 6).Value = 93
    $ws.Cells.Item(14, 6).Value = 222
    $ws.Cells.Item(17, 6).Value = 106
    $ws.Cells.Item(18, 6).Value = 128
    $ws.Cells.Item(19, 6).Value = 3723
    $ws.Cells.Item(20, 6).Value = 6
    $ws.Cells.Item(21, 6).Value = 8
    $ws.Cells.Item(22, 6).Value = 429
    $ws.Cells.Item(23, 6).Value = 341
    $ws.Cells.Item(24, 6).Value = 597
    $ws.Cells.Item(25, 6).Value = 383
    $ws.Cells.Item(26, 6).Value = 350
    $ws.Cells.Item(28, 6).Value = 1517

    # --- 2. Insert a new event row at row 29, pushing the former row 29
    #        ("南昌·代号鸢盛花行only") down to row 30 ---
    $ws.Rows.Item(29).Insert()

    # Copy the formatting of column A's numbering cell down onto the newly
    # inserted row so the style matches the rest of the table.
    $ws.Range("A28").Copy()
    $ws.Range("A29").PasteSpecial(-4122)

    $ws.Cells.Item(29, 1).Value = 28

    # B29 must stay a plain text string ("2024-05-03"), not an auto-converted
    # date serial. Flip the cell to text format before assigning, then strip
    # the temporary number-format override so no extra style id is left on
    # the cell (matches the rest of the B column, which carries no `s` attr).
    $ws.Cells.Item(29, 2).NumberFormat = "@"
    $ws.Cells.Item(29, 2).Value = "2024-05-03"
    $ws.Cells.Item(29, 2).ClearFormats()

    $ws.Cells.Item(29, 3).Value = "新余·LD02国风动漫嘉年华"
    $ws.Cells.Item(29, 4).Value = "劳动北路888号 金联体育篮球馆"
    $ws.Cells.Item(29, 5).Value = "2024.05.03 10:00-05.03 17:00"
    $ws.Cells.Item(29, 6).Value = 5
    $ws.Cells.Item(29, 7).Value = 30
    $ws.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83379"
    $ws.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg"

    # --- 3. Former row 29 is now row 30; update its shifted numbering and
    #        the updated "想去人数" count ---
    $ws.Cells.Item(30, 1).Value = 29
    $ws.Cells.Item(30, 6).Value = 145
}
